# Add a "Dr. " prefix to the Stalford reference line and move the
# "_GoBack" bookmark (which tracks the last edit position) from the end
# of the Kevin Sturm line to the point right before "Harold Stalford".

$d = $word.ActiveDocument

# Locate the insertion point immediately before "Harold Stalford" and
# insert a brand-new run containing "Dr. " there (InsertBefore keeps it
# as its own run rather than merging into the existing "Harold Stalford…"
# run).
$rng = $d.Content
$rng.Find.Execute("Harold Stalford", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0)
$rng.Collapse(1)
$rng.InsertBefore("Dr. ")

# Word keeps exactly one "_GoBack" bookmark, marking the site of the most
# recent edit. Remove it from its old location (end of the Kevin Sturm
# line) …
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# … and re-create it right after the newly inserted "Dr. " run, i.e.
# immediately before "Harold Stalford".
$rng2 = $d.Content
$rng2.Find.Execute("Harold Stalford", $true, $false, $false, $false, $false, `
                    $true, 1, $false, "", 0)
$rng2.Collapse(1)
$d.Bookmarks.Add("_GoBack", $rng2)
